$d = $word.ActiveDocument

# 1. Update the letter date from September 19 to September 21, 2025.
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, `
    $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the single-line mailing address "23305 Eastbrook Ave, Los Altos CA 94024"
#    into two lines: "23305 Eastbrook Ave" and "Los Altos, CA 94024", followed by a
#    new blank line, mirroring the existing paragraph formatting (Arial 11pt,
#    autoSpaceDE/autoSpaceDN off). Using a Find/Replace with paragraph-mark
#    wildcards (^p) keeps the new paragraphs' formatting in sync with the
#    original run instead of leaving stray empty runs behind.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("23305 Eastbrook Ave, Los Altos CA 94024", $false, $false, $false, `
    $false, $false, $true, 1, $false, "23305 Eastbrook Ave^pLos Altos, CA 94024^p", 2) | Out-Null

# 3. Remove the two blank formatting paragraphs that used to sit directly below
#    the "Board of Directors" signature line (a blank No Spacing paragraph and a
#    blank Title-styled paragraph), leaving the following blank Title paragraph
#    untouched.
$count = $d.Paragraphs.Count
$boardIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "Kentfield Pacific Place Board of Directors`r") {
        $boardIndex = $i
        break
    }
}

if ($boardIndex -ge 1) {
    # Delete from the higher index down so the lower index stays valid.
    $firstBlank = $d.Paragraphs.Item($boardIndex + 2)
    $firstBlank.Range.Delete() | Out-Null
    $secondBlank = $d.Paragraphs.Item($boardIndex + 1)
    $secondBlank.Range.Delete() | Out-Null
}
